$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 22:05"

# Update country statistics (values refreshed + table re-sorted descending by Casos totales)
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1809885
$ws.Range("C4").Value = 16355
$ws.Range("D4").Value = 527995
$ws.Range("E4").Value = 1176659
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 689
$ws.Range("H4").Value = 105231

$ws.Range("A72").Value = "Sudan"
$ws.Range("B72").Value = 4800
$ws.Range("C72").Value = 279
$ws.Range("D72").Value = 1272
$ws.Range("E72").Value = 3266
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 29
$ws.Range("H72").Value = 262

$ws.Range("A73").Value = "Guatemala"
$ws.Range("B73").Value = 4607
$ws.Range("C73").Value = 259
$ws.Range("D73").Value = 648
$ws.Range("E73").Value = 3869
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 10
$ws.Range("H73").Value = 90

$ws.Range("A84").Value = "Costa de Marfil"
$ws.Range("B84").Value = 2799
$ws.Range("C84").Value = 49
$ws.Range("D84").Value = 1385
$ws.Range("E84").Value = 1381
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 33

$ws.Range("A99").Value = "Maldivas"
$ws.Range("B99").Value = 1633
$ws.Range("C99").Value = 42
$ws.Range("D99").Value = 386
$ws.Range("E99").Value = 1242
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 5

$ws.Range("A117").Value = "Costa Rica"
$ws.Range("B117").Value = 1047
$ws.Range("C117").Value = 25
$ws.Range("D117").Value = 658
$ws.Range("E117").Value = 379
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 10

$ws.Range("A120").Value = "Republica de Africa Central"
$ws.Range("B120").Value = 962
$ws.Range("C120").Value = 88
$ws.Range("D120").Value = 23
$ws.Range("E120").Value = 938
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 1

$ws.Range("A121").Value = "Niger"
$ws.Range("B121").Value = 955
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 813
$ws.Range("E121").Value = 78
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 64

$ws.Range("A122").Value = "Republica de Chipre"
$ws.Range("B122").Value = 944
$ws.Range("C122").Value = 2
$ws.Range("D122").Value = 790
$ws.Range("E122").Value = 137
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 17

$ws.Range("A139").Value = "Santo Tome y Principe"
$ws.Range("B139").Value = 479
$ws.Range("C139").Value = 16
$ws.Range("D139").Value = 68
$ws.Range("E139").Value = 399
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 12

$ws.Range("A140").Value = "Reunion"
$ws.Range("B140").Value = 471
$ws.Range("C140").Value = 1
$ws.Range("D140").Value = 411
$ws.Range("E140").Value = 59
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 1

$ws.Range("A146").Value = "Ruanda"
$ws.Range("B146").Value = 359
$ws.Range("C146").Value = 4
$ws.Range("D146").Value = 250
$ws.Range("E146").Value = 108
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 1

$ws.Range("A165").Value = "Guyana"
$ws.Range("B165").Value = 150
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 67
$ws.Range("E165").Value = 71
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 1
$ws.Range("H165").Value = 12

$ws.Range("A169").Value = "Libia"
$ws.Range("B169").Value = 130
$ws.Range("C169").Value = 12
$ws.Range("D169").Value = 50
$ws.Range("E169").Value = 75
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 5

$ws.Range("A170").Value = "Camboya"
$ws.Range("B170").Value = 125
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 123
$ws.Range("E170").Value = 2
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 0

$ws.Range("A171").Value = "Siria"
$ws.Range("B171").Value = 122
$ws.Range("C171").Value = 0
$ws.Range("D171").Value = 43
$ws.Range("E171").Value = 75
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 4

$ws.Range("A176").Value = "Monaco"
$ws.Range("B176").Value = 99
$ws.Range("C176").Value = 1
$ws.Range("D176").Value = 90
$ws.Range("E176").Value = 5
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 4

